$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: change existing entry's day/start time, add Category descript + commit hash ---
$ws.Range("A16").Value = 43992
$ws.Range("B16").Value = 0.72986111111111107
$ws.Range("F16").Value = "250059ae9e60f1b754e008bcfc2bdb8743dd5ee1"
$ws.Range("E16").Value = "Category"

# --- Row 17: brand new entry (day/start/end + descript, no commit text) ---
$ws.Range("A17").Value = 43993
$ws.Range("B17").Value = 0.40833333333333338
$ws.Range("C17").Value = 0.50208333333333333
$ws.Range("E17").Value = "Search lawyer by category"

# --- Formatting: column F uses the "commit hash" style (small Consolas font). ---
# Copy that format from an existing styled cell (F2) onto F15:F25 so the
# shared cellXfs entry (s="3") is reused instead of minting new styles.
$ws.Range("F2").Copy()
$ws.Range("F15:F25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection moves to E18 ---
$ws.Range("E18").Select()

Write-Output "done"
